# Adds two new rows to Sheet1 documenting the "Add Digits" and
# "Find the K-th Character in String Game 1" problems.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: Add Digits
# (cell write order mirrors the shared-string insertion order in the target file)
$ws.Range("A23").Value = 258
$ws.Range("B23").Value = "Add Digits"
$ws.Range("C24").Value = "Recursion/Loop/"
$ws.Range("D23").Value = "while in while,"

# Row 24: Find the K-th Character in String Game 1
$ws.Range("A24").Value = 3304
$ws.Range("B24").Value = "Find the K-th Character in String Game 1"
$ws.Range("C23").Value = "Recursion/Loop/Math"
$ws.Range("E23").Value = "1+((num-1) % 9)"
$ws.Range("D24").Value = "Generate word with StringBuilder, return the result[k-1]."

# Match the saved selection/view state (active cell D24, last row added)
[void]$ws.Range("D24").Select()
